# 2018 SAFE pollock parameters: add log_q1_dev/log_q3_dev: columns to "deviates",
# and add a new "main pars" sheet with the remaining model parameters.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("deviates")

# New header cells F1 / G1 on the "deviates" sheet
$ws1.Range("F1").Value = "log_q1_dev"
$ws1.Range("G1").Value = "log_q3_dev:"

# New data columns F2:G50 on the "deviates" sheet
$ws1.Cells.Item(2, 6).Value = 0.00957234715522
$ws1.Cells.Item(2, 7).Value = -0.0547657270549
$ws1.Cells.Item(3, 6).Value = 0.00957234715522
$ws1.Cells.Item(3, 7).Value = -0.0547657270549
$ws1.Cells.Item(4, 6).Value = 0.00957234715522
$ws1.Cells.Item(4, 7).Value = -0.0547657270549
$ws1.Cells.Item(5, 6).Value = 0.00957234715522
$ws1.Cells.Item(5, 7).Value = -0.0547657270549
$ws1.Cells.Item(6, 6).Value = 0.00957234715522
$ws1.Cells.Item(6, 7).Value = -0.0547657270549
$ws1.Cells.Item(7, 6).Value = 0.00957234715522
$ws1.Cells.Item(7, 7).Value = -0.0547657270549
$ws1.Cells.Item(8, 6).Value = 0.00957234715522
$ws1.Cells.Item(8, 7).Value = -0.0547657270549
$ws1.Cells.Item(9, 6).Value = 0.00957234715522
$ws1.Cells.Item(9, 7).Value = -0.0547657270549
$ws1.Cells.Item(10, 6).Value = 0.00957234715522
$ws1.Cells.Item(10, 7).Value = -0.0547657270549
$ws1.Cells.Item(11, 6).Value = 0.00957234715522
$ws1.Cells.Item(11, 7).Value = -0.0547657270549
$ws1.Cells.Item(12, 6).Value = 0.00957234715522
$ws1.Cells.Item(12, 7).Value = -0.0547657270549
$ws1.Cells.Item(13, 6).Value = 0.00957234715522
$ws1.Cells.Item(13, 7).Value = -0.054765727055
$ws1.Cells.Item(14, 6).Value = 0.00957234715522
$ws1.Cells.Item(14, 7).Value = -0.054765727055
$ws1.Cells.Item(15, 6).Value = 0.00957234715522
$ws1.Cells.Item(15, 7).Value = -0.054765727055
$ws1.Cells.Item(16, 6).Value = 0.00957234715522
$ws1.Cells.Item(16, 7).Value = -0.054765727055
$ws1.Cells.Item(17, 6).Value = 0.00957234715523
$ws1.Cells.Item(17, 7).Value = -0.054765727055
$ws1.Cells.Item(18, 6).Value = 0.00957234715523
$ws1.Cells.Item(18, 7).Value = -0.054765727055
$ws1.Cells.Item(19, 6).Value = 0.00957234715523
$ws1.Cells.Item(19, 7).Value = -0.054765727055
$ws1.Cells.Item(20, 6).Value = 0.00957234715523
$ws1.Cells.Item(20, 7).Value = -0.0547657270747
$ws1.Cells.Item(21, 6).Value = 0.00957234715523
$ws1.Cells.Item(21, 7).Value = -0.0564347731608
$ws1.Cells.Item(22, 6).Value = 0.00957234715523
$ws1.Cells.Item(22, 7).Value = -0.0681994095101
$ws1.Cells.Item(23, 6).Value = 0.00957234715523
$ws1.Cells.Item(23, 7).Value = -0.0572535146124
$ws1.Cells.Item(24, 6).Value = 0.00957234715869
$ws1.Cells.Item(24, 7).Value = -0.0487105227672
$ws1.Cells.Item(25, 6).Value = -0.000485253645381
$ws1.Cells.Item(25, 7).Value = -0.0312385125348
$ws1.Cells.Item(26, 6).Value = 0.0147442127183
$ws1.Cells.Item(26, 7).Value = 0.000315801787349
$ws1.Cells.Item(27, 6).Value = 0.050245611232
$ws1.Cells.Item(27, 7).Value = 0.0349970419115
$ws1.Cells.Item(28, 6).Value = 0.0780406634321
$ws1.Cells.Item(28, 7).Value = 0.0704453062732
$ws1.Cells.Item(29, 6).Value = 0.0910747075967
$ws1.Cells.Item(29, 7).Value = 0.113575882524
$ws1.Cells.Item(30, 6).Value = 0.0871057908736
$ws1.Cells.Item(30, 7).Value = 0.151931118133
$ws1.Cells.Item(31, 6).Value = 0.0615399179814
$ws1.Cells.Item(31, 7).Value = 0.193739231289
$ws1.Cells.Item(32, 6).Value = 0.0359740450906
$ws1.Cells.Item(32, 7).Value = 0.234903672576
$ws1.Cells.Item(33, 6).Value = 0.00211334312061
$ws1.Cells.Item(33, 7).Value = 0.261015115713
$ws1.Cells.Item(34, 6).Value = -0.025330361572
$ws1.Cells.Item(34, 7).Value = 0.279873883875
$ws1.Cells.Item(35, 6).Value = -0.054111815666
$ws1.Cells.Item(35, 7).Value = 0.278531548298
$ws1.Cells.Item(36, 6).Value = -0.0878416817143
$ws1.Cells.Item(36, 7).Value = 0.271966918259
$ws1.Cells.Item(37, 6).Value = -0.126021435429
$ws1.Cells.Item(37, 7).Value = 0.257531621786
$ws1.Cells.Item(38, 6).Value = -0.174756581152
$ws1.Cells.Item(38, 7).Value = 0.241716592215
$ws1.Cells.Item(39, 6).Value = -0.221303798385
$ws1.Cells.Item(39, 7).Value = 0.237585021042
$ws1.Cells.Item(40, 6).Value = -0.254016420853
$ws1.Cells.Item(40, 7).Value = 0.229540651369
$ws1.Cells.Item(41, 6).Value = -0.248746533043
$ws1.Cells.Item(41, 7).Value = 0.204837490743
$ws1.Cells.Item(42, 6).Value = -0.206087431651
$ws1.Cells.Item(42, 7).Value = 0.144387888708
$ws1.Cells.Item(43, 6).Value = -0.154825891046
$ws1.Cells.Item(43, 7).Value = 0.0734196098055
$ws1.Cells.Item(44, 6).Value = -0.103564350444
$ws1.Cells.Item(44, 7).Value = -0.00710752601163
$ws1.Cells.Item(45, 6).Value = -0.0127425879673
$ws1.Cells.Item(45, 7).Value = -0.113638658787
$ws1.Cells.Item(46, 6).Value = 0.0676338832167
$ws1.Cells.Item(46, 7).Value = -0.212761569661
$ws1.Cells.Item(47, 6).Value = 0.131210378127
$ws1.Cells.Item(47, 7).Value = -0.315803083946
$ws1.Cells.Item(48, 6).Value = 0.20009805878
$ws1.Cells.Item(48, 7).Value = -0.401240000764
$ws1.Cells.Item(49, 6).Value = 0.289567068718
$ws1.Cells.Item(49, 7).Value = -0.461932164454
$ws1.Cells.Item(50, 6).Value = 0.340322477108
$ws1.Cells.Item(50, 7).Value = -0.465445846033

# Add the new "main pars" worksheet right after "deviates"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "main pars"

$mainParsHeaders = @("log_slp2_srv1", "inf2_srv1", "log_slp1_srv2", "inf1_srv2", "log_slp2_srv2", "inf2_srv2", "log_slp1_srv3", "inf1_srv3", "log_slp1_srv6", "inf1_srv6", "log_slp2_srv6", "inf2_srv6", "log_q1_mean", "log_q2_mean", "log_q3_mean", "log_q4", "log_q5", "log_q6")
for ($i = 0; $i -lt $mainParsHeaders.Count; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $mainParsHeaders[$i]
}

$ws2.Cells.Item(2, 1).Value = 0.296894674756
$ws2.Cells.Item(2, 2).Value = 9.48229975071
$ws2.Cells.Item(2, 3).Value = -0.480950034525
$ws2.Cells.Item(2, 4).Value = 4.11372188292
$ws2.Cells.Item(2, 5).Value = 1.0
$ws2.Cells.Item(2, 6).Value = 20.0
$ws2.Cells.Item(2, 7).Value = 0.354913182705
$ws2.Cells.Item(2, 8).Value = 4.63086247057
$ws2.Cells.Item(2, 9).Value = 4.9
$ws2.Cells.Item(2, 10).Value = 0.5
$ws2.Cells.Item(2, 11).Value = 1.0
$ws2.Cells.Item(2, 12).Value = 20.0
$ws2.Cells.Item(2, 13).Value = -0.494253318292
$ws2.Cells.Item(2, 14).Value = -0.165251283562
$ws2.Cells.Item(2, 15).Value = -0.398566159562
$ws2.Cells.Item(2, 16).Value = -1.09082515064
$ws2.Cells.Item(2, 17).Value = -0.870714003249
$ws2.Cells.Item(2, 18).Value = -0.188670006665

# Approximate the original bestFit column widths for columns A and B
$ws2.Columns.Item(1).ColumnWidth = 13.4
$ws2.Columns.Item(2).ColumnWidth = 11

# Final selections / active sheet, matching the saved view state
$ws1.Range("J5:BF5").Select()
$ws2.Range("M15").Select()
$ws2.Activate()
